$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.617.72'
$ws.Range('E2').Value = '  -3.69%  '
$ws.Range('D3').Value = '2.558.80'
$ws.Range('E3').Value = '  -1.35%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '506.01'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -3.25%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.45'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -6.35%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +0.16%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.554'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -6.80%  '
$ws.Range('D9').Value = '2.561.39'
$ws.Range('E9').Value = '  -1.49%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.18'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -7.84%  '
$ws.Range('E11').Value = '  -3.16%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.332'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -4.18%  '
$ws.Range('E13').Value = '  -0.97%  '
$ws.Range('D14').Value = '3.008.89'
$ws.Range('E14').Value = '  -1.30%  '
$ws.Range('D15').Value = '58.622.35'
$ws.Range('E15').Value = '  -3.72%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '20.57'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -5.01%  '
$ws.Range('E17').Value = '  -4.92%  '
$ws.Range('D18').Value = '2.562.86'
$ws.Range('E18').Value = '  -1.44%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.53'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -4.80%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '333.95'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -5.64%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.08'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -4.57%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.995'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.44%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.94'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -4.49%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '59.37'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -2.33%  '
$ws.Range('E25').Value = '  -4.39%  '
$ws.Range('E26').Value = '  +0.06%  '
$ws.Range('E27').Value = '  -6.43%  '
$ws.Range('D28').Value = '0.0₃0778'
$ws.Range('E28').Value = '  -7.95%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.86'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -7.00%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.00'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +0.03%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.86'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -7.67%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '18.58'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -4.04%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '148.62'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.33%  '
$ws.Range('E34').Value = '  -3.78%  '
$ws.Range('B35').Value = 'SuiNetwork'
$ws.Range('C35').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.910'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -2.34%  '
$ws.Range('B36').Value = 'NEARProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.87'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -6.87%  '
$ws.Range('E37').Value = '  -7.78%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '35.94'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -1.43%  '
$ws.Range('E39').Value = '  -4.66%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.53'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -7.02%  '
$ws.Range('E41').Value = '  -8.34%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '282.63'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -1.90%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.998'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +0.10%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.607'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -2.35%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0980'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -3.41%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0531'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -5.23%  '
$ws.Range('B47').Value = 'WhiteBITCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.33'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.06%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '18.67'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -4.88%  '
$ws.Range('E49').Value = '  -4.53%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '4.52'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -7.39%  '
$ws.Range('D51').Value = '1.912.63'
$ws.Range('E51').Value = '  -2.24%  '
